$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.363.33"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "'2.051.07"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'230.24"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'57.07"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'0.107"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "'14.98"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "'2.356.05"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'20.86"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "'0.757"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "'2.063.05"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'37.224.94"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'69.58"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "'0.0₃0827"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "'227.38"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("D26").Value = "'9.71"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'166.66"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("E28").Value = "  -6.93%  "
$ws.Range("D29").Value = "'19.09"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "'4.55"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").Value = "'0.0619"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").Value = "'4.61"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").Value = "'5.24"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").Value = "'0.0221"
$ws.Range("E40").Value = "  -5.28%  "
$ws.Range("D41").Value = "'1.500.29"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("D42").Value = "'17.11"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0946"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'96.81"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").Value = "'3.96"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "'7.13"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").Value = "'2.92"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").Value = "'2.241.00"
$ws.Range("E51").Value = "  -1.17%  "
